$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 10.9736627452712
$ws.Cells.Item(2, 3).Value = 9.296680535415769
$ws.Cells.Item(2, 4).Value = 3.653592606117823
$ws.Cells.Item(2, 6).Value = 17.59872643132766
$ws.Cells.Item(2, 7).Value = 17.40960513049214
$ws.Cells.Item(2, 8).Value = 11.31221723110832
$ws.Cells.Item(2, 13).Value = 19.96689603056216
$ws.Cells.Item(2, 15).Value = 15.6981614827533

$ws.Cells.Item(3, 2).Value = 10.34268809100116
$ws.Cells.Item(3, 3).Value = 9.05563643969133
$ws.Cells.Item(3, 4).Value = 3.559009926077386
$ws.Cells.Item(3, 6).Value = 17.60423495498907
$ws.Cells.Item(3, 7).Value = 17.38419443477422
$ws.Cells.Item(3, 8).Value = 11.36118454430245
$ws.Cells.Item(3, 13).Value = 19.35867735633029
$ws.Cells.Item(3, 15).Value = 15.76570190521115

$ws.Cells.Item(4, 2).Value = 9.933033659756079
$ws.Cells.Item(4, 3).Value = 8.904159602897293
$ws.Cells.Item(4, 4).Value = 3.499068150005929
$ws.Cells.Item(4, 6).Value = 17.61460389484789
$ws.Cells.Item(4, 7).Value = 17.37900679444185
$ws.Cells.Item(4, 8).Value = 11.39368048097417
$ws.Cells.Item(4, 13).Value = 18.98351723438481
$ws.Cells.Item(4, 15).Value = 15.81228262224021

$ws.Cells.Item(5, 2).Value = 9.760569499825124
$ws.Cells.Item(5, 3).Value = 8.841635043197167
$ws.Cells.Item(5, 4).Value = 3.474194261469712
$ws.Cells.Item(5, 6).Value = 17.62058099442814
$ws.Cells.Item(5, 7).Value = 17.37950296100675
$ws.Cells.Item(5, 8).Value = 11.40753280653787
$ws.Cells.Item(5, 13).Value = 18.83046396632748
$ws.Cells.Item(5, 15).Value = 15.83254289221556

$ws.Cells.Item(6, 2).Value = 9.731600358507709
$ws.Cells.Item(6, 3).Value = 8.831207135963812
$ws.Cells.Item(6, 4).Value = 3.47003765611185
$ws.Cells.Item(6, 6).Value = 17.6216790999104
$ws.Cells.Item(6, 7).Value = 17.37974265155291
$ws.Cells.Item(6, 8).Value = 11.4098697852355
$ws.Cells.Item(6, 13).Value = 18.80504696197679
$ws.Cells.Item(6, 15).Value = 15.83598409339088

$ws.Cells.Item(7, 2).Value = 9.930730028481104
$ws.Cells.Item(7, 3).Value = 8.903319494201519
$ws.Cells.Item(7, 4).Value = 3.498734471513619
$ws.Cells.Item(7, 6).Value = 17.61467741977491
$ws.Cells.Item(7, 7).Value = 17.37900293244252
$ws.Cells.Item(7, 8).Value = 11.39386482992393
$ws.Cells.Item(7, 13).Value = 18.98145345253944
$ws.Cells.Item(7, 15).Value = 15.8125506929967

$ws.Cells.Item(8, 2).Value = 10.76075090311695
$ws.Cells.Item(8, 3).Value = 9.214334492096814
$ws.Cells.Item(8, 4).Value = 3.621381105325589
$ws.Cells.Item(8, 6).Value = 17.59917306954156
$ws.Cells.Item(8, 7).Value = 17.39867820903454
$ws.Cells.Item(8, 8).Value = 11.32859623684467
$ws.Cells.Item(8, 13).Value = 19.75768789472706
$ws.Cells.Item(8, 15).Value = 15.72038464669347

$ws.Cells.Item(9, 2).Value = 12.20992692548608
$ws.Cells.Item(9, 3).Value = 9.793530575333847
$ws.Cells.Item(9, 4).Value = 3.846168979463604
$ws.Cells.Item(9, 6).Value = 17.62435802111663
$ws.Cells.Item(9, 7).Value = 17.52004613545298
$ws.Cells.Item(9, 8).Value = 11.21993002638512
$ws.Cells.Item(9, 13).Value = 21.25571275239735
$ws.Cells.Item(9, 15).Value = 15.580495271353

$ws.Cells.Item(10, 2).Value = 13.16385984303693
$ws.Cells.Item(10, 3).Value = 10.19652466016197
$ws.Cells.Item(10, 4).Value = 4.000644870488932
$ws.Cells.Item(10, 6).Value = 17.67684212485724
$ws.Cells.Item(10, 7).Value = 17.6595426671511
$ws.Cells.Item(10, 8).Value = 11.15193686565066
$ws.Cells.Item(10, 13).Value = 22.32856300691375
$ws.Cells.Item(10, 15).Value = 15.50301294944026

$ws.Cells.Item(11, 2).Value = 13.57354079667325
$ws.Cells.Item(11, 3).Value = 10.3742590541593
$ws.Cells.Item(11, 4).Value = 4.068403248245383
$ws.Cells.Item(11, 6).Value = 17.70807870138606
$ws.Cells.Item(11, 7).Value = 17.73379712328533
$ws.Cells.Item(11, 8).Value = 11.12359040220977
$ws.Cells.Item(11, 13).Value = 22.8082062670447
$ws.Cells.Item(11, 15).Value = 15.47333545959125

$ws.Cells.Item(12, 2).Value = 13.72517498786779
$ws.Cells.Item(12, 3).Value = 10.44070687046231
$ws.Cells.Item(12, 4).Value = 4.093685580921326
$ws.Cells.Item(12, 6).Value = 17.72096118036057
$ws.Cells.Item(12, 7).Value = 17.76344935171436
$ws.Cells.Item(12, 8).Value = 11.11322907864713
$ws.Cells.Item(12, 13).Value = 22.98845109120195
$ws.Cells.Item(12, 15).Value = 15.46290430239257

$ws.Cells.Item(13, 2).Value = 13.69267381138091
$ws.Cells.Item(13, 3).Value = 10.42643501495161
$ws.Cells.Item(13, 4).Value = 4.088257535267003
$ws.Cells.Item(13, 6).Value = 17.71813993140761
$ws.Cells.Item(13, 7).Value = 17.75699535314265
$ws.Cells.Item(13, 8).Value = 11.11544397104523
$ws.Cells.Item(13, 13).Value = 22.94969640788394
$ws.Cells.Item(13, 15).Value = 15.46511484968331

$ws.Cells.Item(14, 2).Value = 13.58608609843354
$ws.Cells.Item(14, 3).Value = 10.37974320192438
$ws.Cells.Item(14, 4).Value = 4.070490862447055
$ws.Cells.Item(14, 6).Value = 17.70911746308846
$ws.Cells.Item(14, 7).Value = 17.73620601749341
$ws.Cells.Item(14, 8).Value = 11.12273048964696
$ws.Cells.Item(14, 13).Value = 22.82306384128571
$ws.Cells.Item(14, 15).Value = 15.47246107310803

$ws.Cells.Item(15, 2).Value = 13.52034153161124
$ws.Cells.Item(15, 3).Value = 10.35103012155247
$ws.Cells.Item(15, 4).Value = 4.059558841419764
$ws.Cells.Item(15, 6).Value = 17.70372801436369
$ws.Cells.Item(15, 7).Value = 17.72367103913757
$ws.Cells.Item(15, 8).Value = 11.12724228945803
$ws.Cells.Item(15, 13).Value = 22.74531232001424
$ws.Cells.Item(15, 15).Value = 15.4770661279766

$ws.Cells.Item(16, 2).Value = 13.13659662943472
$ws.Cells.Item(16, 3).Value = 10.18479221578837
$ws.Cells.Item(16, 4).Value = 3.996164820025617
$ws.Cells.Item(16, 6).Value = 17.67494851181115
$ws.Cells.Item(16, 7).Value = 17.65490571528282
$ws.Cells.Item(16, 8).Value = 11.15384148041535
$ws.Cells.Item(16, 13).Value = 22.2970325437194
$ws.Cells.Item(16, 15).Value = 15.50506505194338

$ws.Cells.Item(17, 2).Value = 12.89495470194098
$ws.Cells.Item(17, 3).Value = 10.08134002988551
$ws.Cells.Item(17, 4).Value = 3.956619841106567
$ws.Cells.Item(17, 6).Value = 17.6591758874212
$ws.Cells.Item(17, 7).Value = 17.61547379022151
$ws.Cells.Item(17, 8).Value = 11.17082195933178
$ws.Cells.Item(17, 13).Value = 22.01974352737308
$ws.Cells.Item(17, 15).Value = 15.52367247964435

$ws.Cells.Item(18, 2).Value = 12.7536848661358
$ws.Cells.Item(18, 3).Value = 10.02131376578147
$ws.Cells.Item(18, 4).Value = 3.933638827479688
$ws.Cells.Item(18, 6).Value = 17.65079717052194
$ws.Cells.Item(18, 7).Value = 17.59381127061784
$ws.Cells.Item(18, 8).Value = 11.18083181923735
$ws.Cells.Item(18, 13).Value = 21.85947286993717
$ws.Cells.Item(18, 15).Value = 15.53489886894307

$ws.Cells.Item(19, 2).Value = 12.70546158058402
$ws.Cells.Item(19, 3).Value = 10.00090170915039
$ws.Cells.Item(19, 4).Value = 3.92581783108372
$ws.Cells.Item(19, 6).Value = 17.64807946425061
$ws.Cells.Item(19, 7).Value = 17.58665199965307
$ws.Cells.Item(19, 8).Value = 11.18426270544636
$ws.Cells.Item(19, 13).Value = 21.8050794851921
$ws.Cells.Item(19, 15).Value = 15.5387897168835

$ws.Cells.Item(20, 2).Value = 12.92091437503131
$ws.Cells.Item(20, 3).Value = 10.09240724797048
$ws.Cells.Item(20, 4).Value = 3.960853986999737
$ws.Cells.Item(20, 6).Value = 17.66078318640165
$ws.Cells.Item(20, 7).Value = 17.61956618182708
$ws.Cells.Item(20, 8).Value = 11.16898918517997
$ws.Cells.Item(20, 13).Value = 22.04934362608216
$ws.Cells.Item(20, 15).Value = 15.52163742103782

$ws.Cells.Item(21, 2).Value = 13.61748864753335
$ws.Cells.Item(21, 3).Value = 10.39348135456754
$ws.Cells.Item(21, 4).Value = 4.075719688700855
$ws.Cells.Item(21, 6).Value = 17.71173902567649
$ws.Cells.Item(21, 7).Value = 17.74227089931351
$ws.Cells.Item(21, 8).Value = 11.12058013266531
$ws.Cells.Item(21, 13).Value = 22.86029782842141
$ws.Cells.Item(21, 15).Value = 15.47028135822849

$ws.Cells.Item(22, 2).Value = 14.05231872688731
$ws.Cells.Item(22, 3).Value = 10.58523831765008
$ws.Cells.Item(22, 4).Value = 4.148591843456315
$ws.Cells.Item(22, 6).Value = 17.75118086407615
$ws.Cells.Item(22, 7).Value = 17.83139297904923
$ws.Cells.Item(22, 8).Value = 11.09111597483639
$ws.Cells.Item(22, 13).Value = 23.38215670305294
$ws.Cells.Item(22, 15).Value = 15.44142401440018

$ws.Cells.Item(23, 2).Value = 13.82211342727715
$ws.Cells.Item(23, 3).Value = 10.48336875300637
$ws.Cells.Item(23, 4).Value = 4.109904371122726
$ws.Cells.Item(23, 6).Value = 17.72957029726409
$ws.Cells.Item(23, 7).Value = 17.78301741449909
$ws.Cells.Item(23, 8).Value = 11.10664220455261
$ws.Cells.Item(23, 13).Value = 23.10443026294243
$ws.Cells.Item(23, 15).Value = 15.45639309749545

$ws.Cells.Item(24, 2).Value = 12.90918531511975
$ws.Cells.Item(24, 3).Value = 10.08740546845987
$ws.Cells.Item(24, 4).Value = 3.958940494978111
$ws.Cells.Item(24, 6).Value = 17.66005437957434
$ws.Cells.Item(24, 7).Value = 17.61771287231964
$ws.Cells.Item(24, 8).Value = 11.16981701138927
$ws.Cells.Item(24, 13).Value = 22.03596406627888
$ws.Cells.Item(24, 15).Value = 15.52255582462493

$ws.Cells.Item(25, 2).Value = 11.8372356522409
$ws.Cells.Item(25, 3).Value = 9.640545068381009
$ws.Cells.Item(25, 4).Value = 3.787160045113734
$ws.Cells.Item(25, 6).Value = 17.61157697127513
$ws.Cells.Item(25, 7).Value = 17.47834516751758
$ws.Cells.Item(25, 8).Value = 11.2472517573267
$ws.Cells.Item(25, 13).Value = 20.85445576512474
$ws.Cells.Item(25, 15).Value = 15.61392450226412
